# Adapt column header formatting to respective input file names.
# Renames the "_old" / "_new" header-name suffixes to the respective
# format-version suffixes ("_FV2310" / "_FV2404"), wraps the sheet's used
# range in an Excel Table (ListObject), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1, columns A:U) from "<name>_old"/"<name>_new"
#    to "<name>_FV2310"/"<name>_FV2404". "diff" (column K) stays untouched.
$headerRange = $ws.Range("A1:U1")
for ($i = 1; $i -le $headerRange.Columns.Count; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $text = $cell.Value2
    if ($text -ne $null) {
        if ($text.EndsWith("_old")) {
            $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2310"
        } elseif ($text.EndsWith("_new")) {
            $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2404"
        }
    }
}

# 2. Turn the A1:U61 range into an Excel Table so the new headers are also
#    used as the table's column names / structured-reference headers.
$tableRange = $ws.Range("A1:U61")
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# 3. Freeze the header row (split below row 1).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
